$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1047
$ws1.Range("F3").Value = 338
$ws1.Range("F4").Value = 1451
$ws1.Range("F5").Value = 8647
$ws1.Range("F8").Value = 634
$ws1.Range("F11").Value = 3499
$ws1.Range("F13").Value = 355
$ws1.Range("F15").Value = 1085
$ws1.Range("F19").Value = 183
$ws1.Range("F20").Value = 2223
$ws1.Range("F21").Value = 40

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1047
$ws4.Range("F3").Value = 338
$ws4.Range("F4").Value = 1451
$ws4.Range("F5").Value = 8647
$ws4.Range("F8").Value = 634
$ws4.Range("F11").Value = 3499
$ws4.Range("F13").Value = 355
$ws4.Range("F15").Value = 1085
$ws4.Range("F19").Value = 183
$ws4.Range("F20").Value = 2223
$ws4.Range("F22").Value = 40
